$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "69.857.63"
$ws.Range("E2").Value = "  -1.07%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.533.03"
$ws.Range("E3").Value = "  -1.21%  "

$ws.Range("E4").Value = "  -0.07%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "607.67"
$ws.Range("E5").Value = "  +3.13%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "185.07"
$ws.Range("E6").Value = "  -1.38%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.530.94"
$ws.Range("E7").Value = "  -0.96%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.613"
$ws.Range("E8").Value = "  -1.67%  "

$ws.Range("E9").Value = "  -0.05%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.212"
$ws.Range("E10").Value = "  +5.99%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.640"
$ws.Range("E11").Value = "  -1.51%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "53.54"
$ws.Range("E12").Value = "  -2.48%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000307"
$ws.Range("E13").Value = "  -0.56%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.44"
$ws.Range("E14").Value = "  -1.82%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.091.52"
$ws.Range("E15").Value = "  -1.28%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "69.893.10"
$ws.Range("E16").Value = "  -0.96%  "

$ws.Range("B17").Value = "Uniswap"
$ws.Range("C17").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "12.60"
$ws.Range("E17").Value = "  +1.00%  "

$ws.Range("B18").Value = "Chainlink"
$ws.Range("C18").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "18.84"
$ws.Range("E18").Value = "  -3.26%  "

$ws.Range("B19").Value = "WrappedEther"
$ws.Range("C19").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.531.02"
$ws.Range("E19").Value = "  -1.50%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "574.13"
$ws.Range("E20").Value = "  +2.81%  "

$ws.Range("E21").Value = "  -0.21%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.989"
$ws.Range("E22").Value = "  -3.27%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "17.35"
$ws.Range("E23").Value = "  -2.63%  "

$ws.Range("E24").Value = "  -0.40%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.86"
$ws.Range("E25").Value = "  -1.95%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "93.87"
$ws.Range("E26").Value = "  -2.08%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.97"
$ws.Range("E27").Value = "  -0.76%  "

$ws.Range("E28").Value = "  -5.10%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.34"
$ws.Range("E29").Value = "  +1.70%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "32.02"
$ws.Range("E30").Value = "  -0.87%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.98"
$ws.Range("E31").Value = "  -5.09%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "12.14"
$ws.Range("E32").Value = "  -3.40%  "

$ws.Range("E33").Value = "  -1.76%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "63.15"
$ws.Range("E34").Value = "  -3.07%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.28"
$ws.Range("E35").Value = "  +0.05%  "

$ws.Range("E36").Value = "  +16.60%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "532.83"
$ws.Range("E37").Value = "  -4.33%  "

$ws.Range("E38").Value = "  -3.71%  "

$ws.Range("E39").Value = "  +0.08%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "37.00"
$ws.Range("E40").Value = "  -3.31%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.543.36"
$ws.Range("E41").Value = "  +5.69%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0₃0777"
$ws.Range("E42").Value = "  +0.36%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.51"
$ws.Range("E43").Value = "  +3.78%  "

$ws.Range("E44").Value = "  +0.54%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0451"
$ws.Range("E45").Value = "  +1.05%  "

$ws.Range("B46").Value = "ThetaToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.91"
$ws.Range("E46").Value = "  -2.70%  "

$ws.Range("B47").Value = "ApeXProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.42"
$ws.Range("E47").Value = "  -4.05%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.06"
$ws.Range("E49").Value = "  -4.03%  "

$ws.Range("E50").Value = "  +0.24%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.42"
$ws.Range("E51").Value = "  -3.91%  "
